$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.0
$ws.Range("B2").Value = -0.07513316298689968
$ws.Range("C2").Value = -0.0
$ws.Range("D2").Value = 0.224556992242958
$ws.Range("E2").Value = 0.008134815887270613
$ws.Range("G2").Value = 0.0
$ws.Range("I2").Value = -0.0
$ws.Range("J2").Value = -0.0
$ws.Range("K2").Value = -0.003065432192118038
$ws.Range("L2").Value = -0.0
$ws.Range("M2").Value = 0.2181395890659293
$ws.Range("N2").Value = -0.0001727015592207062
$ws.Range("S2").Value = 0.0
$ws.Range("T2").Value = -0.1034251572261255
$ws.Range("V2").Value = 0.01448340324324757
$ws.Range("W2").Value = -0.0237150189021015
$ws.Range("Z2").Value = -0.0
$ws.Range("AB2").Value = 0.0
$ws.Range("AC2").Value = -0.05533071069590436
$ws.Range("AD2").Value = 0.0
$ws.Range("AE2").Value = -0.01713963613298442
$ws.Range("AF2").Value = -0.01542473018745792
$ws.Range("AG2").Value = -0.0
$ws.Range("AI2").Value = -0.0
$ws.Range("AJ2").Value = 0.0
$ws.Range("AK2").Value = -0.0
$ws.Range("AL2").Value = -0.02499626156132146
$ws.Range("AM2").Value = 0.0
$ws.Range("AN2").Value = 0.02433931786697912
$ws.Range("AO2").Value = 0.07335784195797591
$ws.Range("AQ2").Value = 0.0
$ws.Range("AR2").Value = -0.0
$ws.Range("AT2").Value = 0.0
$ws.Range("AU2").Value = -0.09747175097281018
$ws.Range("AW2").Value = 0.07048834296269527
$ws.Range("AX2").Value = 0.01136902336123896
$ws.Range("AY2").Value = -0.0
$ws.Range("BC2").Value = -0.0
$ws.Range("BD2").Value = -0.02865649821570666
$ws.Range("BF2").Value = 0.1033052057149321
$ws.Range("BG2").Value = 0.02854527687070273
$ws.Range("BI2").Value = -0.0
$ws.Range("BJ2").Value = -0.0
$ws.Range("BL2").Value = 0.0
$ws.Range("BM2").Value = 0.02405885502493103
$ws.Range("BO2").Value = -0.03544817232043314
$ws.Range("BP2").Value = -0.06990616414759654
$ws.Range("BU2").Value = 0.0
$ws.Range("BV2").Value = -0.04038669193535912
$ws.Range("BW2").Value = 0.0
$ws.Range("BX2").Value = 0.02249102905515309
$ws.Range("BY2").Value = -0.007291196476252951
$ws.Range("BZ2").Value = -0.0
$ws.Range("CB2").Value = 0.0
$ws.Range("CD2").Value = -0.0
$ws.Range("CE2").Value = 0.03040352147995489
$ws.Range("CG2").Value = -0.03267992352547659
$ws.Range("CH2").Value = 0.01224640038650156
$ws.Range("CJ2").Value = -0.0
$ws.Range("CM2").Value = -0.0
$ws.Range("CN2").Value = -0.009288310456953388
$ws.Range("CO2").Value = -0.0
$ws.Range("CP2").Value = 0.01161331693495279
$ws.Range("CQ2").Value = 0.04950123115183769
$ws.Range("CT2").Value = 0.0
$ws.Range("CU2").Value = -0.0
$ws.Range("CV2").Value = -0.0
$ws.Range("CW2").Value = 0.04530621560229153
$ws.Range("CY2").Value = -0.03669277790065874
$ws.Range("CZ2").Value = 0.007527598420719271
$ws.Range("DE2").Value = -0.0
$ws.Range("DF2").Value = 0.0332190769767132
$ws.Range("DH2").Value = 0.03651431811600674
$ws.Range("DI2").Value = 0.04671352483237037
$ws.Range("DJ2").Value = 0.0
$ws.Range("DK2").Value = -0.0
$ws.Range("DL2").Value = -0.0
$ws.Range("DN2").Value = 0.0
$ws.Range("DO2").Value = -0.02344918865452078
$ws.Range("DP2").Value = -0.0
$ws.Range("DQ2").Value = 0.03415356360395851
$ws.Range("DR2").Value = -0.03563786140452697
$ws.Range("DS2").Value = -0.0
$ws.Range("DW2").Value = 0.0
$ws.Range("DX2").Value = -0.06148322476389054
$ws.Range("DY2").Value = -0.0
$ws.Range("DZ2").Value = -0.01348045115148365
$ws.Range("EA2").Value = -0.03216953766349701
$ws.Range("EB2").Value = 0.0
$ws.Range("EF2").Value = -0.0
$ws.Range("EG2").Value = 0.04248559568891849
$ws.Range("EI2").Value = 0.06517213132250246
$ws.Range("EJ2").Value = -0.02494558276727713
$ws.Range("EO2").Value = 0.0
$ws.Range("EP2").Value = 0.04440061486491791
$ws.Range("EQ2").Value = 0.0
$ws.Range("ER2").Value = -0.04343541494434484
$ws.Range("ES2").Value = 0.04427119156428989
$ws.Range("ET2").Value = 0.0
$ws.Range("EU2").Value = -0.0
$ws.Range("EV2").Value = 0.0
$ws.Range("EX2").Value = 0.0
$ws.Range("EY2").Value = 0.04175983852563794
$ws.Range("EZ2").Value = 0.0
$ws.Range("FA2").Value = -0.02837712393765177
$ws.Range("FB2").Value = 0.01515671653930182
$ws.Range("FD2").Value = -0.0
$ws.Range("FG2").Value = -0.0
$ws.Range("FH2").Value = -0.001727621869630965
$ws.Range("FJ2").Value = -0.004742519334957517
$ws.Range("FK2").Value = -0.004203330212255114
$ws.Range("FL2").Value = -0.0
$ws.Range("FN2").Value = -0.0
$ws.Range("FP2").Value = -0.0
$ws.Range("FQ2").Value = -0.0112934643681217
$ws.Range("FR2").Value = -0.0
$ws.Range("FS2").Value = -0.01983029870617169
$ws.Range("FT2").Value = -0.001772087722976833
$ws.Range("FV2").Value = -0.0
$ws.Range("FY2").Value = 0.0
$ws.Range("FZ2").Value = -0.04271389418570019
$ws.Range("GB2").Value = 0.02963769923369615
$ws.Range("GD2").Value = 0.0
$ws.Range("GE2").Value = -0.0
